# Update Pertanggal 17 Januari 2023 17:38 WIB
# Adds new "Project ..." roles to the TblAppObject_UserRole SQL generator sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: named parameter binding (-Row 5 -RoleName "x") on custom functions is
# unreliable on this host, so this helper is always called positionally.
function Set-RoleRow {
    param($Row, $RoleName)
    $ws.Cells.Item($Row, 2).Value = $RoleName
    $ws.Cells.Item($Row, 3).Formula = '=IF(EXACT(B' + $Row + ',""),"",CONCATENATE("PERFORM ""SchSysConfig"".""Func_TblAppObject_UserRole_SET""(varSystemLoginSession, null, null, null, varInstitutionBranchID, ''",B' + $Row + ',"'');"))'
}

# --- Insert the first two new rows (Project General Manager / Project Management Officer)
# right before the existing "Finance, Accounting, & Tax General Manager" row (row 17).
$ws.Rows("17:18").Insert()

Set-RoleRow 17 "Project General Manager"
$ws.Cells.Item(17, 4).Value = 95000000000013

Set-RoleRow 18 "Project Management Officer"
$ws.Cells.Item(18, 4).Value = 95000000000014

# Renumber the D column (ID) for the rows that shifted down by 2.
$ws.Cells.Item(19, 4).Value = 95000000000015
$ws.Cells.Item(20, 4).Value = 95000000000016
$ws.Cells.Item(21, 4).Value = 95000000000017

# --- Insert nine new rows: eight Project roles (Project Manager .. Project Worker
# Staff) plus the trailing blank separator row, right before the existing
# "Finance Manager" row (now at row 23 after the first insert).
$ws.Rows("23:31").Insert()

Set-RoleRow 23 "Project Manager"
$ws.Cells.Item(23, 4).Value = 95000000000018

Set-RoleRow 24 "Site Manager"
$ws.Cells.Item(24, 4).Value = 95000000000019

Set-RoleRow 25 "Project Controller Supervisor"
$ws.Cells.Item(25, 4).Value = 95000000000020

Set-RoleRow 26 "Project Controller Staff"
$ws.Cells.Item(26, 4).Value = 95000000000021

Set-RoleRow 27 "Project Administrator Supervisor"
$ws.Cells.Item(27, 4).Value = 95000000000022

Set-RoleRow 28 "Project Administrator Staff"
$ws.Cells.Item(28, 4).Value = 95000000000023

Set-RoleRow 29 "Project Worker Supervisor"
$ws.Cells.Item(29, 4).Value = 95000000000024

Set-RoleRow 30 "Project Worker Staff"
$ws.Cells.Item(30, 4).Value = 95000000000025

# Row 31 stays blank (separator), same as the other section breaks.

# Renumber the D column (ID) for the remaining rows that shifted down by 11 total.
$ws.Cells.Item(32, 4).Value = 95000000000026
$ws.Cells.Item(33, 4).Value = 95000000000027
$ws.Cells.Item(34, 4).Value = 95000000000028
$ws.Cells.Item(35, 4).Value = 95000000000029
$ws.Cells.Item(36, 4).Value = 95000000000030
$ws.Cells.Item(38, 4).Value = 95000000000031
$ws.Cells.Item(39, 4).Value = 95000000000032
$ws.Cells.Item(40, 4).Value = 95000000000033
$ws.Cells.Item(42, 4).Value = 95000000000034
$ws.Cells.Item(43, 4).Value = 95000000000035
$ws.Cells.Item(44, 4).Value = 95000000000036
$ws.Cells.Item(45, 4).Value = 95000000000037
$ws.Cells.Item(46, 4).Value = 95000000000038

# Restore the selection shown by the author at save time.
$ws.Range("F41").Select()
